$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows -------------------------------------------------

# Row 2: product name changed, barcode now filled in
$ws.Range("A2").Value = "PRODUTO 12 EXEMPLO TESTE 12 - INICIO"
$ws.Range("B2").Value = 6989652331394

# Row 3: product name changed
$ws.Range("A3").Value = "CABO PRODUTO 13 EXEMPLO"

# Row 8: barcode corrected
$ws.Range("B8").Value = 6989652331394

# --- Append new product rows -------------------------------------------

# Row 11 mirrors the formatting of row 9 (empty barcode style)
$ws.Range("A9:F9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "CABO PRODUTO 9 EXEMPLO 09"
$ws.Range("C11").Value = 277
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 5

# Row 12 mirrors the formatting of row 10 (filled barcode style)
$ws.Range("A10:F10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "CARREGADOR PRODUTO 10 EXEMPLO 010"
$ws.Range("B12").Value = 10663543819178
$ws.Range("C12").Value = 360.84
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 8

# Row 13 mirrors the formatting of row 10 (filled barcode style)
$ws.Range("A10:F10").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "PRODUTO PRODUTO 11 EXEMPLO 011 - FINAL"
$ws.Range("B13").Value = 10663543819465
$ws.Range("C13").Value = 360.84
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 8

$excel.CutCopyMode = 0

# --- Selection mirrors the saved workbook state -------------------------
$null = $ws.Range("B2").Select()
